# Weekly fruit/vegetable price update: two new records inserted at the top
# of the data block (rows 49-50), pushing the existing rows down by 2
# (old row 49 -> new row 51, ... old row 87 -> new row 89).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 49 (twice, inserting at
# the same index each time pushes the previously-inserted row further down).
$ws.Rows.Item(49).Insert()
$ws.Rows.Item(49).Insert()

# New row 49
$ws.Cells.Item(49, 1).Value = 10
$ws.Cells.Item(49, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(49, 3).Value = "La Araucanía"
$ws.Cells.Item(49, 4).Value = 44923
$ws.Cells.Item(49, 5).Value = 9
$ws.Cells.Item(49, 6).Value = 100112030
$ws.Cells.Item(49, 7).Value = "Poroto granado"
$ws.Cells.Item(49, 8).Value = "Sin especificar"
$ws.Cells.Item(49, 9).Value = "Primera"
$ws.Cells.Item(49, 10).Value = 35
$ws.Cells.Item(49, 11).Value = 43000
$ws.Cells.Item(49, 12).Value = 43000
$ws.Cells.Item(49, 13).Value = 43000
$ws.Cells.Item(49, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(49, 15).Value = "Región del Maule"
$ws.Cells.Item(49, 16).Value = 1720
$ws.Cells.Item(49, 17).Value = 25
$ws.Cells.Item(49, 18).Value = "Hortaliza"

# New row 50
$ws.Cells.Item(50, 1).Value = 10
$ws.Cells.Item(50, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(50, 3).Value = "La Araucanía"
$ws.Cells.Item(50, 4).Value = 44923
$ws.Cells.Item(50, 5).Value = 9
$ws.Cells.Item(50, 6).Value = 100112030
$ws.Cells.Item(50, 7).Value = "Poroto granado"
$ws.Cells.Item(50, 8).Value = "Sin especificar"
$ws.Cells.Item(50, 9).Value = "Primera"
$ws.Cells.Item(50, 10).Value = 35
$ws.Cells.Item(50, 11).Value = 43000
$ws.Cells.Item(50, 12).Value = 43000
$ws.Cells.Item(50, 13).Value = 43000
$ws.Cells.Item(50, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(50, 15).Value = "Región del Maule"
$ws.Cells.Item(50, 16).Value = 1720
$ws.Cells.Item(50, 17).Value = 25
$ws.Cells.Item(50, 18).Value = "Hortaliza"
